{"js": "// Apply the edit described by the diff to the first paragraph of the body:\n//  1. Replace the placeholder id text\n//     **ID__AFFARS_5334_topic_4__ID**  ->  **ID__AFFARS_SUBPART_5334_70__ID**\n//  2. Drop the now-unwanted trailing \" \" run that followed the id text.\n//  3. Change the paragraph's left indent from 120 twips (6pt) to 225 twips (11.25pt).\n//  4. Add a paragraph border (top/left/bottom/right) that reserves 5 twips of\n//     spacing (\"w:space\") on every side - Word.ParagraphBorder (Office.js) has\n//     no \"space\"/\"distance\" property, so that attribute is applied via a\n//     scoped insertOoxml() on the same paragraph once the text/indent edits\n//     above are already in place.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// --- 1. Replace the id placeholder text -----------------------------------\nconst idRange = firstParagraph.getRange();\nconst idMatches = idRange.search(\"**ID__AFFARS_5334_topic_4__ID**\", { matchCase: true });\nidMatches.load(\"items\");\nawait context.sync();\n\nidMatches.items[0].insertText(\"**ID__AFFARS_SUBPART_5334_70__ID**\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 2. Remove the trailing space run that used to follow the id text -----\nconst trailingRange = firstParagraph.getRange();\nconst spaceMatches = trailingRange.search(\" \", { matchCase: true });\nspaceMatches.load(\"items\");\nawait context.sync();\n\nif (spaceMatches.items.length > 0) {\n  spaceMatches.items[0].delete();\n  await context.sync();\n}\n\n// --- 3. Update the left indent (120 -> 225 twips = 6pt -> 11.25pt) --------\nfirstParagraph.leftIndent = 11.25;\nawait context.sync();\n\n// --- 4. Add the pBdr (top/left/bottom/right, 5-twip space, no line style) -\n// Word.ParagraphBorder only exposes color/type/width/visible - there is no\n// way to set the OOXML \"w:space\" attribute through the regular Office.js\n// object model, so the border is spliced in directly via insertOoxml on the\n// (already edited) paragraph.\nconst borderOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pBdr>\n                <w:top w:space=\"5\"/>\n                <w:left w:space=\"5\"/>\n                <w:bottom w:space=\"5\"/>\n                <w:right w:space=\"5\"/>\n              </w:pBdr>\n              <w:spacing w:after=\"0\"/>\n              <w:ind w:left=\"225\"/>\n              <w:jc w:val=\"left\"/>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/>\n                <w:b w:val=\"false\"/>\n                <w:i w:val=\"false\"/>\n                <w:color w:val=\"000000\"/>\n                <w:sz w:val=\"22\"/>\n              </w:rPr>\n              <w:t>**ID__AFFARS_SUBPART_5334_70__ID**</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nfirstParagraph.insertOoxml(borderOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Apply the edit described by the diff to the first paragraph of the document:\n#  1. Replace the placeholder id text\n#     **ID__AFFARS_5334_topic_4__ID**  ->  **ID__AFFARS_SUBPART_5334_70__ID**\n#  2. Drop the now-unwanted trailing \" \" run that followed the id text.\n#  3. Change the paragraph's left indent from 120 twips (6pt) to 225 twips (11.25pt).\n#  4. Add a paragraph border (top/left/bottom/right) that reserves 5pt of\n#     spacing (\"w:space\") on every side, with no explicit line style/color.\n\n$d = $word.ActiveDocument\n$p1 = $d.Paragraphs(1)\n\n# --- 1. Replace the id placeholder text, keeping the trailing space run ---\n$idLen = \"**ID__AFFARS_5334_topic_4__ID**\".Length\n$idRange = $d.Range(0, $idLen)\n$idRange.Text = \"**ID__AFFARS_SUBPART_5334_70__ID**\"\n\n# --- 2. Remove the trailing space run that used to follow the id text -----\n$newIdLen = \"**ID__AFFARS_SUBPART_5334_70__ID**\".Length\n$trailingRange = $d.Range($newIdLen, $newIdLen + 1)\n$trailingRange.Delete()\n\n# --- 3. Update the left indent (120 -> 225 twips = 6pt -> 11.25pt) --------\n$p1.LeftIndent = 11.25\n\n# --- 4. Add the paragraph border (5pt space on every side) ----------------\n$borders = $p1.Borders\n$borders.DistanceFromTop = 5\n$borders.DistanceFromBottom = 5\n$borders.DistanceFromLeft = 5\n$borders.DistanceFromRight = 5\n"}
